# Update odds values in Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 (Ind. Rivadavia - Union de Santa Fe)
$ws.Range("G6").Value = 2.63
$ws.Range("I6").Value = 3.1
$ws.Range("M6").Value = 1.17
$ws.Range("N6").Value = 5
$ws.Range("O6").Value = 1.67
$ws.Range("P6").Value = 2.1
$ws.Range("S6").Value = 5.4
$ws.Range("T6").Value = 1.16
$ws.Range("U6").Value = 7
$ws.Range("V6").Value = 1.1
$ws.Range("W6").Value = 1.67
$ws.Range("AA6").Value = 5.5
$ws.Range("AD6").Value = 26
$ws.Range("AG6").Value = 5

# Row 7 (Talleres Cordoba - Rosario Central)
$ws.Range("G7").Value = 1.95
$ws.Range("I7").Value = 4.5
$ws.Range("J7").Value = 2.75
$ws.Range("W7").Value = 1.57
$ws.Range("AB7").Value = 8
$ws.Range("AD7").Value = 17
$ws.Range("AR7").Value = 1.93
$ws.Range("AS7").Value = 1.93

# Row 15 (Deportes Tolima - Once Caldas)
$ws.Range("G15").Value = 1.65
$ws.Range("H15").Value = 3.75
$ws.Range("I15").Value = 5.25
$ws.Range("K15").Value = 2.25
$ws.Range("L15").Value = 5.5
$ws.Range("Q15").Value = 2.03
$ws.Range("R15").Value = 1.83
$ws.Range("S15").Value = 2.85
$ws.Range("T15").Value = 1.41
$ws.Range("U15").Value = 3.5
$ws.Range("V15").Value = 1.29
$ws.Range("W15").Value = 1.36
$ws.Range("X15").Value = 3
$ws.Range("Y15").Value = 1.91
$ws.Range("Z15").Value = 1.8
$ws.Range("AB15").Value = 7.5
$ws.Range("AD15").Value = 12
$ws.Range("AE15").Value = 15
$ws.Range("AH15").Value = 7
$ws.Range("AI15").Value = 17
$ws.Range("AM15").Value = 26
$ws.Range("AN15").Value = 17
$ws.Range("AQ15").Value = 41
